# Generate Report for Handoff
# The c834d112-8d1f-43ed-b454-6a8a13c240f5 file moves from
# "Handed back: in sync with en-US" to "Ready for handoff" with an
# updated handoff datetime, on the Overview sheet as well as the
# zh-cn and de-de detail sheets.

$wb = $excel.ActiveWorkbook

$status = "Ready for handoff"

# --- Overview sheet (row 3 = c834d112-...md) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $status
$wsOverview.Range("C3").Value = $status
$wsOverview.Range("D3").Value = "2016-03-25 08:19:12"

# --- zh-cn sheet (row 3 = c834d112-...md) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $status
$wsZhCn.Range("E3").Value = "2016-03-25 08:19:02"

# --- de-de sheet (row 3 = c834d112-...md) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $status
$wsDeDe.Range("E3").Value = "2016-03-25 08:19:12"
